# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale sheets, flips the
# Overview/locale "Status" text from "Ready for handoff" to
# "Handed back: in sync with en-US", and widens a few columns so the new
# long file-name values are readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d5041af5504156aff321d7d7a63bc74a97df00c/e2e/86bd36d9-8eff-480d-8f23-7f4ce7cbbbb3.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d5041af5504156aff321d7d7a63bc74a97df00c/e2e/ca70f728-2df9-4f16-9b6d-b0c6eee1d14c.md"
$mdName1 = "86bd36d9-8eff-480d-8f23-7f4ce7cbbbb3.md"
$mdName2 = "ca70f728-2df9-4f16-9b6d-b0c6eee1d14c.md"

# ---------------------------------------------------------------------
# Overview sheet: flip the zh-cn / de-de status cells, widen E:F columns
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText
$ov.Columns.Item(5).ColumnWidth = 29.16666666666667
$ov.Columns.Item(6).ColumnWidth = 29.16666666666667

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Range("I2").Value = $mdName1
$zh.Range("J2").Value = "86bd36d9-8eff-480d-8f23-7f4ce7cbbbb3.ec9ac459316ccc900346498703b6e455856a9d22.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-17 20:48:51"

$zh.Range("I3").Value = $mdName2
$zh.Range("J3").Value = "ca70f728-2df9-4f16-9b6d-b0c6eee1d14c.6aedd07443e925261bd6869c30415e4900bcb9e2.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-17 20:48:51"

$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl1, "", "", $mdName1)
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl2, "", "", $mdName2)

$zh.Range("I2").Font.Underline = 2
$zh.Range("I2").Font.Color = 15570276
$zh.Range("I3").Font.Underline = 2
$zh.Range("I3").Font.Color = 15570276

$zh.Columns.Item(3).ColumnWidth = 29.16666666666667
$zh.Columns.Item(9).ColumnWidth = 39.16666666666667
$zh.Columns.Item(10).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Range("I2").Value = $mdName1
$de.Range("J2").Value = "86bd36d9-8eff-480d-8f23-7f4ce7cbbbb3.ec9ac459316ccc900346498703b6e455856a9d22.de-de.xlf"
$de.Range("K2").Value = "2016-08-17 20:48:59"

$de.Range("I3").Value = $mdName2
$de.Range("J3").Value = "ca70f728-2df9-4f16-9b6d-b0c6eee1d14c.6aedd07443e925261bd6869c30415e4900bcb9e2.de-de.xlf"
$de.Range("K3").Value = "2016-08-17 20:48:59"

$de.Hyperlinks.Add($de.Range("I2"), $mdUrl1, "", "", $mdName1)
$de.Hyperlinks.Add($de.Range("I3"), $mdUrl2, "", "", $mdName2)

$de.Range("I2").Font.Underline = 2
$de.Range("I2").Font.Color = 15570276
$de.Range("I3").Font.Underline = 2
$de.Range("I3").Font.Color = 15570276

$de.Columns.Item(3).ColumnWidth = 29.16666666666667
$de.Columns.Item(9).ColumnWidth = 39.16666666666667
$de.Columns.Item(10).ColumnWidth = 39.16666666666667
